# Applies the updated currentAveragePrice / Leve cost-profit recalculation values
# captured for the Anima_Profits workbook scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1884.3636
$ws.Range("J107").Value = 1227.6
$ws.Range("L107").Value = 1227.6
$ws.Range("N107").Value = -5067.6
$ws.Range("H116").Value = 3400.353
$ws.Range("I116").Value = 3066.6667
$ws.Range("J116").Value = 3775.75
$ws.Range("K116").Value = 3066.6667
$ws.Range("L116").Value = 3775.75
$ws.Range("M116").Value = 375.3332999999998
$ws.Range("N116").Value = -10659.75
$ws.Range("H125").Value = 1263.6666
$ws.Range("I125").Value = 610.8333
$ws.Range("J125").Value = 2569.3333
$ws.Range("K125").Value = 5497.4997
$ws.Range("L125").Value = 23123.9997
$ws.Range("M125").Value = -3037.4997
$ws.Range("N125").Value = -28043.9997
$ws.Range("H138").Value = 1972.7467
$ws.Range("I138").Value = 2013.7391
$ws.Range("J138").Value = 1954.6154
$ws.Range("K138").Value = 6041.2173
$ws.Range("L138").Value = 5863.8462
$ws.Range("M138").Value = -901.2173000000003
$ws.Range("N138").Value = -16143.8462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6795.4
$ws.Range("I132").Value = 8252
$ws.Range("J132").Value = 4610.5
$ws.Range("K132").Value = 24756
$ws.Range("L132").Value = 13831.5
$ws.Range("M132").Value = -22226
$ws.Range("N132").Value = -18891.5
$ws.Range("H134").Value = 55428.5
$ws.Range("J134").Value = 55428.5
$ws.Range("L134").Value = 55428.5
$ws.Range("N134").Value = -65568.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 71431370
$ws.Range("I86").Value = 90912010
$ws.Range("J86").Value = 2333.3333
$ws.Range("K86").Value = 90912010
$ws.Range("L86").Value = 2333.3333
$ws.Range("M86").Value = -90910887
$ws.Range("N86").Value = -4579.3333
$ws.Range("H89").Value = 71431370
$ws.Range("I89").Value = 90912010
$ws.Range("J89").Value = 2333.3333
$ws.Range("K89").Value = 454560050
$ws.Range("L89").Value = 11666.6665
$ws.Range("M89").Value = -454554434
$ws.Range("N89").Value = -22898.6665
$ws.Range("H134").Value = 5237.5
$ws.Range("I134").Value = 6725
$ws.Range("K134").Value = 20175
$ws.Range("M134").Value = -17640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9999
$ws.Range("J51").Value = 9999
$ws.Range("L51").Value = 9999
$ws.Range("N51").Value = -11471
$ws.Range("H59").Value = 31999
$ws.Range("J59").Value = 31999
$ws.Range("L59").Value = 31999
$ws.Range("N59").Value = -34289
$ws.Range("H61").Value = 9999
$ws.Range("J61").Value = 9999
$ws.Range("L61").Value = 9999
$ws.Range("N61").Value = -10695
$ws.Range("H68").Value = 22682.75
$ws.Range("J68").Value = 23817.545
$ws.Range("L68").Value = 23817.545
$ws.Range("N68").Value = -25315.545
$ws.Range("H71").Value = 22682.75
$ws.Range("J71").Value = 23817.545
$ws.Range("L71").Value = 71452.63499999999
$ws.Range("N71").Value = -78940.63499999999
$ws.Range("H74").Value = 25599.5
$ws.Range("J74").Value = 29499.375
$ws.Range("L74").Value = 29499.375
$ws.Range("N74").Value = -31247.375
$ws.Range("H77").Value = 25599.5
$ws.Range("J77").Value = 29499.375
$ws.Range("L77").Value = 88498.125
$ws.Range("N77").Value = -97234.125
$ws.Range("H132").Value = 11113448
$ws.Range("I132").Value = 1515.2
$ws.Range("J132").Value = 33337314
$ws.Range("K132").Value = 4545.6
$ws.Range("L132").Value = 100011942
$ws.Range("M132").Value = -2015.6
$ws.Range("N132").Value = -100017002
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 6000
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 18000
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -15465
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 150000
$ws.Range("J37").Value = 150000
$ws.Range("L37").Value = 450000
$ws.Range("N37").Value = -450224
$ws.Range("H58").Value = 1624.7368
$ws.Range("I58").Value = 820
$ws.Range("J58").Value = 1669.4445
$ws.Range("K58").Value = 2460
$ws.Range("L58").Value = 5008.333500000001
$ws.Range("M58").Value = -2332
$ws.Range("N58").Value = -5264.333500000001
$ws.Range("H70").Value = 3506
$ws.Range("I70").Value = 3012
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 9036
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -8721
$ws.Range("N70").Value = -12630
$ws.Range("H73").Value = 3506
$ws.Range("I73").Value = 3012
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 9036
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -7944
$ws.Range("N73").Value = -14184
$ws.Range("H131").Value = 933.4516
$ws.Range("I131").Value = 228.42857
$ws.Range("J131").Value = 1139.0834
$ws.Range("K131").Value = 685.28571
$ws.Range("L131").Value = 3417.2502
$ws.Range("M131").Value = 4354.71429
$ws.Range("N131").Value = -13497.2502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 212120.8
$ws.Range("I11").Value = 500300
$ws.Range("K11").Value = 500300
$ws.Range("M11").Value = -500161
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").Value = $null
$ws.Range("H80").Value = 1178496.8
$ws.Range("I80").Value = 1504162.4
$ws.Range("J80").Value = 201500
$ws.Range("K80").Value = 1504162.4
$ws.Range("L80").Value = 201500
$ws.Range("M80").Value = -1503164.4
$ws.Range("N80").Value = -203496
$ws.Range("H83").Value = 1178496.8
$ws.Range("I83").Value = 1504162.4
$ws.Range("J83").Value = 201500
$ws.Range("K83").Value = 7520812
$ws.Range("L83").Value = 1007500
$ws.Range("M83").Value = -7515820
$ws.Range("N83").Value = -1017484
$ws.Range("H93").Value = 52499.168
$ws.Range("J93").Value = 52499.168
$ws.Range("L93").Value = 52499.168
$ws.Range("N93").Value = -56243.168
$ws.Range("H97").Value = 2813.875
$ws.Range("I97").Value = 2400
$ws.Range("J97").Value = 3503.6667
$ws.Range("K97").Value = 2400
$ws.Range("L97").Value = 3503.6667
$ws.Range("M97").Value = -1904
$ws.Range("N97").Value = -4495.6667
$ws.Range("H132").Value = 3199.9333
$ws.Range("I132").Value = 3066.6667
$ws.Range("J132").Value = 3733
$ws.Range("K132").Value = 9200.000100000001
$ws.Range("L132").Value = 11199
$ws.Range("M132").Value = -6670.000100000001
$ws.Range("N132").Value = -16259

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8452.532999999999
$ws.Range("I22").Value = 1261.2858
$ws.Range("J22").Value = 14744.875
$ws.Range("K22").Value = 1261.2858
$ws.Range("L22").Value = 14744.875
$ws.Range("M22").Value = -966.2858000000001
$ws.Range("N22").Value = -15334.875
$ws.Range("H27").Value = 8452.532999999999
$ws.Range("I27").Value = 1261.2858
$ws.Range("J27").Value = 14744.875
$ws.Range("K27").Value = 1261.2858
$ws.Range("L27").Value = 14744.875
$ws.Range("M27").Value = -1154.2858
$ws.Range("N27").Value = -14958.875
$ws.Range("H100").Value = 3156.6
$ws.Range("I100").Value = 2901.5
$ws.Range("J100").Value = 3326.6667
$ws.Range("K100").Value = 2901.5
$ws.Range("L100").Value = 3326.6667
$ws.Range("M100").Value = -2360.5
$ws.Range("N100").Value = -4408.6667
$ws.Range("H132").Value = 4237.5947
$ws.Range("I132").Value = 3258.3684
$ws.Range("J132").Value = 5271.222
$ws.Range("K132").Value = 9775.1052
$ws.Range("L132").Value = 15813.666
$ws.Range("M132").Value = -7245.1052
$ws.Range("N132").Value = -20873.666
$ws.Range("H136").Value = 10419236
$ws.Range("I136").Value = 3525.25
$ws.Range("J136").Value = 20834946
$ws.Range("K136").Value = 10575.75
$ws.Range("L136").Value = 62504838
$ws.Range("M136").Value = -8025.75
$ws.Range("N136").Value = -62509938

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11500
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 13666.667
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 13666.667
$ws.Range("M41").Value = -4610
$ws.Range("N41").Value = -14446.667
$ws.Range("H132").Value = 5379113.5
$ws.Range("I132").Value = 5229.2856
$ws.Range("J132").Value = 6946496.5
$ws.Range("K132").Value = 15687.8568
$ws.Range("L132").Value = 20839489.5
$ws.Range("M132").Value = -13157.8568
$ws.Range("N132").Value = -20844549.5
$ws.Range("H136").Value = 3676.7273
$ws.Range("I136").Value = 3409.5334
$ws.Range("K136").Value = 10228.6002
$ws.Range("M136").Value = -7678.600199999999
